$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("Salary") - shifts old "Apply Here" link column from D to E
$ws.Columns.Item(4).Insert()

# Remove the now-stale rows 11-16 (old data had 15 rows, new data only has 9)
$ws.Range("A11:E16").EntireRow.Delete()

# Header row
$ws.Range("A1").Value = 'Title'
$ws.Range("B1").Value = 'Company'
$ws.Range("C1").Value = 'Location'
$ws.Range("D1").Value = 'Salary'
$ws.Range("E1").Value = 'Apply Here'

# Row 2
$ws.Range("A2").Value = 'Associate Software Developer Apprentice'
$ws.Range("B2").Value = 'Pearson'
$ws.Range("C2").Value = 'Remote in London'
$ws.Range("D2").Value = '£22,000 - £25,000 a year'
$ws.Hyperlinks.Add($ws.Range("E2"), 'https://uk.indeed.com/rc/clk?jk=002bf4846fcb7d28&fccid=915b1c0ee87e5e8a&vjs=3', "", "", 'test/rc/clk?jk=002bf4846fcb7d28&fccid=915b1c0ee87e5e8a&vjs=3')
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("A3").Value = 'Software Junior Developer Apprentice'
$ws.Range("B3").Value = 'QA Apprenticeships'
$ws.Range("C3").Value = 'London EC2M'
$ws.Range("D3").Value = '£25,000 - £30,000 a year'
$ws.Hyperlinks.Add($ws.Range("E3"), 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0CXy9_N1TLC2ejqyTX_V4eQ4PmQ2039NyFlBjhJ2joCdiLJ2guK6p6gDxKa49r6zxjCTAdQ1uM6iXGzvNU-Yc7zz3hcZmLelAG_vJQZlpFH00o1s04wldLCcuG3xoLJUNvg5I_mfxXCc_PV8WtCkXK1fUnA5Ex_hentUerjmxq2xmbvt24wPEOhAt-4w0krKuY-5OC1YS1GG8ekoDS0pUvI9LB72Jbk-px0W1IHgYf2nnyKhygMNemCiNhxr6fwaA69LYKM45kVJZ45QUPOiGNyc0P3Bxd6RuqbJbi1M5HxoOAYwevuidoX8c1_R9ukVmVc48n3TDL0IMNmdc4N9vufm8IkpwGJZqwDrWj10JKh0p8tBBcpC2W2oaGSpfeBrtAhQTcwhgXFO7H2sMiaJ1SHql9Cg8qrSW8wW3OuRbyC9D1g_j8Y575urpomtH8gy-_h4m2_9wVcvLHsCr4CADU40sJAalV6IQXDY7FBh-hzlTmEi_TLXsHTyH5gt7ZglxClV7w7I8qPjXY0KCA00rNqdbwUZLvf0lypRVaTQDJk9tmrkLwVsDtJ&xkcb=SoDK-_M3c1NbXUwHsp0LbzkdCdPP&p=4&fvj=0&vjs=3', "", "", 'test/pagead/clk?mo=r&ad=-6NYlbfkN0CXy9_N1TLC2ejqyTX_V4eQ4PmQ2039NyFlBjhJ2joCdiLJ2guK6p6gDxKa49r6zxjCTAdQ1uM6iXGzvNU-Yc7zz3hcZmLelAG_vJQZlpFH00o1s04wldLCcuG3xoLJUNvg5I_mfxXCc_PV8WtCkXK1fUnA5Ex_hentUerjmxq2xmbvt24wPEOhAt-4w0krKuY-5OC1YS1GG8ekoDS0pUvI9LB72Jbk-px0W1IHgYf2nnyKhygMNemCiNhxr6fwaA69LYKM45kVJZ45QUPOiGNyc0P3Bxd6RuqbJbi1M5HxoOAYwevuidoX8c1_R9ukVmVc48n3TDL0IMNmdc4N9vufm8IkpwGJZqwDrWj10JKh0p8tBBcpC2W2oaGSpfeBrtAhQTcwhgXFO7H2sMiaJ1SHql9Cg8qrSW8wW3OuRbyC9D1g_j8Y575urpomtH8gy-_h4m2_9wVcvLHsCr4CADU40sJAalV6IQXDY7FBh-hzlTmEi_TLXsHTyH5gt7ZglxClV7w7I8qPjXY0KCA00rNqdbwUZLvf0lypRVaTQDJk9tmrkLwVsDtJ&xkcb=SoDK-_M3c1NbXUwHsp0LbzkdCdPP&p=4&fvj=0&vjs=3')
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("A4").Value = 'We are recruiting an Apprentice Software Developer'
$ws.Range("B4").Value = 'thecitysecret'
$ws.Range("C4").Value = 'Richmond'
$ws.Range("D4").Value = '£4.40 an hour'
$ws.Hyperlinks.Add($ws.Range("E4"), 'https://uk.indeed.com/rc/clk?jk=63779fc7e4e1e5f2&fccid=905a967fd25ae49d&vjs=3', "", "", 'test/rc/clk?jk=63779fc7e4e1e5f2&fccid=905a967fd25ae49d&vjs=3')
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("A5").Value = 'Software Engineer Apprentice'
$ws.Range("B5").Value = 'QA Apprenticeships'
$ws.Range("C5").Value = 'London TW8'
$ws.Range("D5").Value = '£20,000 a year'
$ws.Hyperlinks.Add($ws.Range("E5"), 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0CXy9_N1TLC2ejqyTX_V4eQ4PmQ2039NyFlBjhJ2joCdiLJ2guK6p6gDxKa49r6zxjSg5XGI2hT26G5hch9mNuzoY5w7goGqmcNfC2DDEHXlsp8_Lo7x8TSxtXFBDW0VOlKt7Kf0mszN0utQPq0jsLsqgwtOPc52bi3BJrDqmGQyGnXwDlHBKqty8jiZnhuGBkJITRleNl6DU2WS2ckhTixJD0_K6DPaJyR6Y6Keb8EofWnAWcdbVKaOXK9XQ5quwPXVSK3GV3Bo_p-EjGFj77zsf0ZM4_Htg-mIqZBAKkPpl4BJMHcZXoUBKmwo15zCZjhrvw7JIxNcs_ARWI5bYBGp4s6BRM_E1CDdwSRBT6qslTWx2X7Knr4QqIPtOUavsz6kttv0XST2a2PEZ_Yt8mrhTRQjJEkaiXt85aUpzredDZs_g17PfWGJp0TgwNtSxugOifJf6ggvPznCl92EVSx2EtkwDzpDo2w2TGtHIiuFKrp3nx3iwaSNasT7oK7EmtX0f--CBVH15ztrOylAxJIYj-qCV2IXX7iddRhbvPdSTAdTpNOwBXp&xkcb=SoAn-_M3c1NbXXQHsp0KbzkdCdPP&p=6&fvj=0&vjs=3', "", "", 'test/pagead/clk?mo=r&ad=-6NYlbfkN0CXy9_N1TLC2ejqyTX_V4eQ4PmQ2039NyFlBjhJ2joCdiLJ2guK6p6gDxKa49r6zxjSg5XGI2hT26G5hch9mNuzoY5w7goGqmcNfC2DDEHXlsp8_Lo7x8TSxtXFBDW0VOlKt7Kf0mszN0utQPq0jsLsqgwtOPc52bi3BJrDqmGQyGnXwDlHBKqty8jiZnhuGBkJITRleNl6DU2WS2ckhTixJD0_K6DPaJyR6Y6Keb8EofWnAWcdbVKaOXK9XQ5quwPXVSK3GV3Bo_p-EjGFj77zsf0ZM4_Htg-mIqZBAKkPpl4BJMHcZXoUBKmwo15zCZjhrvw7JIxNcs_ARWI5bYBGp4s6BRM_E1CDdwSRBT6qslTWx2X7Knr4QqIPtOUavsz6kttv0XST2a2PEZ_Yt8mrhTRQjJEkaiXt85aUpzredDZs_g17PfWGJp0TgwNtSxugOifJf6ggvPznCl92EVSx2EtkwDzpDo2w2TGtHIiuFKrp3nx3iwaSNasT7oK7EmtX0f--CBVH15ztrOylAxJIYj-qCV2IXX7iddRhbvPdSTAdTpNOwBXp&xkcb=SoAn-_M3c1NbXXQHsp0KbzkdCdPP&p=6&fvj=0&vjs=3')
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("A6").Value = 'Robotics Developer'
$ws.Range("B6").Value = 'Barnet and Chase Farm Hospitals'
$ws.Range("C6").Value = 'London'
$ws.Range("D6").Value = '£45,024 - £50,806 a year'
$ws.Hyperlinks.Add($ws.Range("E6"), 'https://uk.indeed.com/rc/clk?jk=e897606392e83984&fccid=7691ebb71b24124c&vjs=3', "", "", 'test/rc/clk?jk=e897606392e83984&fccid=7691ebb71b24124c&vjs=3')
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("A7").Value = 'Software Developer Veteran Training Programme'
$ws.Range("B7").Value = 'SaluteMyJob'
$ws.Range("C7").Value = 'Woking'
$ws.Range("D7").Value = '£30,000 - £45,000 a year'
$ws.Hyperlinks.Add($ws.Range("E7"), 'https://uk.indeed.com/company/SaluteMyJob/jobs/Software-Developer-Veteran-Training-Programme-2e5bfa9ee0f81f49?fccid=111f6b9664e56375&vjs=3', "", "", 'test/company/SaluteMyJob/jobs/Software-Developer-Veteran-Training-Programme-2e5bfa9ee0f81f49?fccid=111f6b9664e56375&vjs=3')
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("A8").Value = 'Software Support Engineer Degree Apprentice'
$ws.Range("B8").Value = 'QA Apprenticeships'
$ws.Range("C8").Value = 'Weybridge KT13'
$ws.Range("D8").Value = '£25,000 - £28,000 a year'
$ws.Hyperlinks.Add($ws.Range("E8"), 'https://uk.indeed.com/pagead/clk?mo=r&ad=-6NYlbfkN0CXy9_N1TLC2ejqyTX_V4eQ4PmQ2039NyFlBjhJ2joCdiLJ2guK6p6gDxKa49r6zxgiO854Ph9UtTNrnGuT_m6j4f4lvnPNlmnEtP0fthAefN3dVdqyYs1SAo_xq3EkMME-v2zAmZpURi2zEIn-tZnXw1K2iiPBVzrXrg_dOjKLAPekfv3DUM3izRKf3xk_KWgLNqczCkYih765V8AAhW8NCVCd-Ubq_Jftim8l-jq2VPQqoiBO9PWF7KGDhYj5_YIa6vzEY2Ff0R_8kTdX0x3DvKYQwHKkjBqg6CCmMUTQnZeF8xS4HcF42-YBOK44ALfKgZsR4M-mappq4kfG_9q3sZFaVth0BJWuG67-TpcUprGoSu6T69r4XRo-B_kiAyaD6_uqPj-vF43MycEEms-spPSb_hzwl97JTDnQChqOIsrcu7cmBep2kphAtu_z6Q4UhTEzqcFESzBo5XklwhVp1QGLGe2m15ce0kM7---tfJFKiewP-M3hRdzX92NHW2hs4juWcYJEtuzBQPfk3MuxwzAVMei40KQjfsaFctCahpCL4C0RVPsd&xkcb=SoBY-_M3c1NbXWQHsp0LbzkdCdPP&p=10&fvj=0&vjs=3', "", "", 'test/pagead/clk?mo=r&ad=-6NYlbfkN0CXy9_N1TLC2ejqyTX_V4eQ4PmQ2039NyFlBjhJ2joCdiLJ2guK6p6gDxKa49r6zxgiO854Ph9UtTNrnGuT_m6j4f4lvnPNlmnEtP0fthAefN3dVdqyYs1SAo_xq3EkMME-v2zAmZpURi2zEIn-tZnXw1K2iiPBVzrXrg_dOjKLAPekfv3DUM3izRKf3xk_KWgLNqczCkYih765V8AAhW8NCVCd-Ubq_Jftim8l-jq2VPQqoiBO9PWF7KGDhYj5_YIa6vzEY2Ff0R_8kTdX0x3DvKYQwHKkjBqg6CCmMUTQnZeF8xS4HcF42-YBOK44ALfKgZsR4M-mappq4kfG_9q3sZFaVth0BJWuG67-TpcUprGoSu6T69r4XRo-B_kiAyaD6_uqPj-vF43MycEEms-spPSb_hzwl97JTDnQChqOIsrcu7cmBep2kphAtu_z6Q4UhTEzqcFESzBo5XklwhVp1QGLGe2m15ce0kM7---tfJFKiewP-M3hRdzX92NHW2hs4juWcYJEtuzBQPfk3MuxwzAVMei40KQjfsaFctCahpCL4C0RVPsd&xkcb=SoBY-_M3c1NbXWQHsp0LbzkdCdPP&p=10&fvj=0&vjs=3')
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("A9").Value = 'Programmer and Developer Apprentice'
$ws.Range("B9").Value = 'Penknife Integrated Marketing'
$ws.Range("C9").Value = 'Watford'
$ws.Range("D9").Value = '£10,000 - £12,000 a year'
$ws.Hyperlinks.Add($ws.Range("E9"), 'https://uk.indeed.com/rc/clk?jk=7b8644857d4b5a7d&fccid=1bbce0be4428fd2a&vjs=3', "", "", 'test/rc/clk?jk=7b8644857d4b5a7d&fccid=1bbce0be4428fd2a&vjs=3')
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("A10").Value = 'Software Development Teacher'
$ws.Range("B10").Value = 'Big Creative Training Ltd'
$ws.Range("C10").Value = 'London'
$ws.Range("D10").Value = '£27,000 - £33,000 a year'
$ws.Hyperlinks.Add($ws.Range("E10"), 'https://uk.indeed.com/company/Big-Creative-Education---Apprenticeships/jobs/Software-Development-Teacher-34283d1df4002ff6?fccid=50ad3812a1dd6750&vjs=3', "", "", 'test/company/Big-Creative-Education---Apprenticeships/jobs/Software-Development-Teacher-34283d1df4002ff6?fccid=50ad3812a1dd6750&vjs=3')
$ws.Range("E10").Style = "Normal"